# Re-cache the auto-updating "datetimeFigureOut" date fields that live on
# the notes master, the slide master and every slide layout's Date
# Placeholder. PowerPoint refreshes the cached <a:t> text of these fields
# on save; this commit's capture day moved from 17.01.21 to 05.03.21.
#
# Also makes sure the presentation-level slide-guide list extension is
# present (PowerPoint writes an empty p15:sldGuideLst the first time the
# guides feature is touched during a session).

$newDate = "05.03.21"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$p = $ppt.ActivePresentation

# Notes master "Date Placeholder" field.
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide master "Date Placeholder" field.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's "Date Placeholder" field (layouts without one are
# simply skipped by Update-DatePlaceholder's guard).
foreach ($layout in $p.SlideMaster.CustomLayouts) {
    Update-DatePlaceholder $layout.Shapes
}

# Touch the presentation-level guide list so PowerPoint's empty
# <p15:sldGuideLst/> extension is (re)written on save, matching a
# session where the Guides feature was opened.
$guides = $p.Guides
if ($guides -ne $null) {
    $g = $guides.Add(1, 1800)
    if ($g -ne $null) {
        $g.Delete()
    }
}
